$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restore full numeric precision for existing computed cells (cosmetic re-save) ---
$ws.Range("Q2").Value2 = 1.7532931884002589
$ws.Range("S2").Value2 = 5.7735849056603774
$ws.Range("T2").Value2 = 6.8837209302325579
$ws.Range("U2").Value2 = 1.9291593386304979
$ws.Range("V2").Value2 = 170.04614843888859
$ws.Range("N3").Value2 = 0.97727272727272729
$ws.Range("S3").Value2 = 4.5199999999999996
$ws.Range("T3").Value2 = 5.9142857142857146
$ws.Range("U3").Value2 = 1.7773707317759559
$ws.Range("V3").Value2 = 109.79202438784159
$ws.Range("Z3").Value2 = 1.0056657223796031
$ws.Range("N4").Value2 = 0.92592592592592593
$ws.Range("R4").Value2 = 84.675296923255345
$ws.Range("S4").Value2 = 3.7931034482758621
$ws.Range("T4").Value2 = 5.5454545454545459
$ws.Range("U4").Value2 = 1.7129785913749409
$ws.Range("V4").Value2 = 93.471706484626964
$ws.Range("Z4").Value2 = 1.0236686390532539
$ws.Range("Q5").Value2 = 1.7268015729532831
$ws.Range("R5").Value2 = 153.47951663347601
$ws.Range("S5").Value2 = 5.6226415094339623
$ws.Range("T5").Value2 = 5.6226415094339623
$ws.Range("U5").Value2 = 1.7268015729532831
$ws.Range("V5").Value2 = 153.47951663347601
$ws.Range("N6").Value2 = 0.99636363636363634
$ws.Range("Q6").Value2 = 2.1086170308137668
$ws.Range("S6").Value2 = 8.2368421052631575
$ws.Range("T6").Value2 = 8.8285714285714292
$ws.Range("U6").Value2 = 2.1779932154083319
$ws.Range("V6").Value2 = 197.77023746070839
$ws.Range("N7").Value2 = 0.87878787878787878
$ws.Range("Q7").Value2 = 1.5064386729619541
$ws.Range("R7").Value2 = 94.197382370788162
$ws.Range("S7").Value2 = 4.5106382978723403
$ws.Range("T7").Value2 = 6.3703703703703702
$ws.Range("U7").Value2 = 1.8516576108091241
$ws.Range("V7").Value2 = 95.005244508153652
$ws.Range("Z7").Value2 = 1.0344827586206899
$ws.Range("R8").Value2 = 91.874137616907788
$ws.Range("S8").Value2 = 3.3012048192771091
$ws.Range("T8").Value2 = 3.3012048192771091
$ws.Range("V8").Value2 = 91.874137616907788
$ws.Range("N9").Value2 = 0.84090909090909094
$ws.Range("R9").Value2 = 67.400815917567527
$ws.Range("U9").Value2 = 1.4971087274601811
$ws.Range("V9").Value2 = 63.092520721274219
$ws.Range("Z9").Value2 = 1.0363321799307961
$ws.Range("R10").Value2 = 166.56550359384281
$ws.Range("S10").Value2 = 5.0735294117647056
$ws.Range("T10").Value2 = 5.0735294117647056
$ws.Range("V10").Value2 = 166.56550359384281
$ws.Range("R11").Value2 = 227.02298809493189
$ws.Range("S11").Value2 = 6.6833333333333336
$ws.Range("T11").Value2 = 6.6833333333333336
$ws.Range("V11").Value2 = 227.02298809493189

# --- New summary row 12: average of column J ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$jFont = $ws.Range("J12").Font
$jFont.Bold = $true
$jFont.Size = 11

# --- New summary rows 14-17: labels + aggregate formulas ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$summaryRange = $ws.Range("B14:B17")
$summaryFont = $summaryRange.Font
$summaryFont.Bold = $true
$summaryFont.Size = 12
$summaryRange.VerticalAlignment = -4108

# --- Row heights for the new bold/larger-font summary rows ---
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# --- Selection + page setup to mirror the saved view state ---
$ws.Range("A14:B17").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
